# Update the cryptos list (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) hold plain text in this sheet, e.g.
# "69.870.41" or "1.00" -- values that Excel would otherwise auto-convert
# to a number. Prefix with a leading apostrophe (exactly like typing
# '69.87... into a cell) so the value stays literal text.
function Set-Text($row, $col, $text) {
    if ($null -eq $text) { return }
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

function Set-Row($row, $price, $volume) {
    Set-Text $row 4 $price
    Set-Text $row 5 $volume
}

Set-Row 2  "69.727.61"  "  +1.96%  "
Set-Row 3  "3.500.46"   "  +0.94%  "
Set-Row 4  $null         "  -0.06%  "
Set-Row 5  "604.97"     "  +2.92%  "
Set-Row 6  "174.54"     "  +4.02%  "
Set-Row 7  $null         "  +0.49%  "
Set-Row 8  "3.496.01"   "  +1.03%  "
Set-Row 9  $null         "  +0.00%  "
Set-Row 10 "0.191"      "  +0.81%  "
Set-Row 11 "7.28"       "  +8.00%  "
Set-Row 12 "0.581"      "  +1.91%  "
Set-Row 13 "46.24"      "  -0.32%  "
Set-Row 14 $null         "  -0.35%  "
Set-Row 15 "4.067.20"   "  +0.83%  "
Set-Row 16 "8.28"       "  -0.12%  "
Set-Row 17 "609.46"     "  -0.68%  "
Set-Row 18 "3.505.42"   $null
Set-Row 19 "69.769.67"  "  +1.89%  "
Set-Row 21 "17.15"      "  -0.07%  "
Set-Row 22 "0.871"      "  +0.26%  "
Set-Row 23 "9.06"       "  -18.19%  "
Set-Row 24 "15.43"      "  -1.70%  "
Set-Row 25 "95.66"      "  +0.16%  "
Set-Row 26 "3.72"       "  -1.26%  "
Set-Row 27 "0.999"      "  -0.10%  "
Set-Row 28 "2.56"       "  -1.30%  "
Set-Row 29 "33.92"      "  +3.98%  "
Set-Row 30 "8.95"       "  -1.20%  "
Set-Row 31 "688.84"     "  +20.82%  "

# Rows 32 and 33 swap places: Stacks <-> Filecoin.
$ws.Cells.Item(32, 2).Value = "Filecoin"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-Text 32 4 "8.10"
Set-Text 32 5 "  -3.21%  "

$ws.Cells.Item(33, 2).Value = "Stacks"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-Text 33 4 "2.97"
Set-Text 33 5 "  -3.07%  "

Set-Row 34 "6.90"       "  +1.86%  "
Set-Row 35 $null         "  -2.62%  "
Set-Row 36 "0.0997"     "  -0.78%  "
Set-Row 37 "3.53"       "  +1.57%  "
Set-Row 38 $null         "  +0.41%  "
Set-Row 39 $null         "  +8.60%  "

# Rows 40 and 41 swap places: FirstDigitalUSD <-> OKB.
$ws.Cells.Item(40, 2).Value = "OKB"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-Text 40 4 "56.55"
Set-Text 40 5 "  -0.58%  "

$ws.Cells.Item(41, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-Text 41 4 "1.00"
Set-Text 41 5 "  +0.01%  "

Set-Row 42 $null         "  +4.75%  "
Set-Row 43 "3.315.86"   "  -1.94%  "
Set-Row 44 $null         "  -2.97%  "
Set-Row 45 "2.93"       "  +4.73%  "
Set-Row 46 "32.24"      "  -0.70%  "
Set-Row 47 $null         "  +0.02%  "
Set-Row 48 "2.54"       "  +0.41%  "
Set-Row 49 $null         "  +1.33%  "
Set-Row 50 "133.37"     "  +0.76%  "
Set-Row 51 $null         "  -0.05%  "
